# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
# Updates Price (D) and Volume(1h) (E) columns for the cryptos list;
# rows 43/44 (Fetch.AI / Kaspa) swap order with updated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.796.04"
$ws.Cells.Item(2, 5).Value = "  -4.10%  "
$ws.Cells.Item(3, 4).Value = "3.384.85"
$ws.Cells.Item(3, 5).Value = "  -4.59%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).Value = "'562.56"
$ws.Cells.Item(5, 5).Value = "  -4.11%  "
$ws.Cells.Item(6, 4).Value = "'184.70"
$ws.Cells.Item(6, 5).Value = "  -7.13%  "
$ws.Cells.Item(7, 5).Value = "  -1.97%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).Value = "3.375.73"
$ws.Cells.Item(9, 5).Value = "  -4.47%  "
$ws.Cells.Item(10, 4).Value = "'0.189"
$ws.Cells.Item(10, 5).Value = "  -8.42%  "
$ws.Cells.Item(11, 5).Value = "  -4.66%  "
$ws.Cells.Item(12, 4).Value = "'48.40"
$ws.Cells.Item(12, 5).Value = "  -7.24%  "
$ws.Cells.Item(13, 5).Value = "  -6.46%  "
$ws.Cells.Item(14, 5).Value = "  -5.86%  "
$ws.Cells.Item(15, 4).Value = "3.922.54"
$ws.Cells.Item(15, 5).Value = "  -4.52%  "
$ws.Cells.Item(16, 4).Value = "'611.07"
$ws.Cells.Item(16, 5).Value = "  -11.28%  "
$ws.Cells.Item(17, 4).Value = "'18.44"
$ws.Cells.Item(17, 5).Value = "  -0.81%  "
$ws.Cells.Item(18, 4).Value = "66.757.81"
$ws.Cells.Item(18, 5).Value = "  -4.22%  "
$ws.Cells.Item(19, 4).Value = "3.379.83"
$ws.Cells.Item(19, 5).Value = "  -5.01%  "
$ws.Cells.Item(20, 5).Value = "  -2.90%  "
$ws.Cells.Item(21, 4).Value = "'11.66"
$ws.Cells.Item(21, 5).Value = "  -6.59%  "
$ws.Cells.Item(22, 5).Value = "  -5.10%  "
$ws.Cells.Item(23, 4).Value = "'17.06"
$ws.Cells.Item(23, 5).Value = "  -4.84%  "
$ws.Cells.Item(24, 4).Value = "'5.23"
$ws.Cells.Item(24, 5).Value = "  +0.13%  "
$ws.Cells.Item(25, 4).Value = "'99.27"
$ws.Cells.Item(25, 5).Value = "  -8.57%  "
$ws.Cells.Item(26, 5).Value = "  -6.65%  "
$ws.Cells.Item(27, 4).Value = "'6.02"
$ws.Cells.Item(27, 5).Value = "  +0.41%  "
$ws.Cells.Item(28, 5).Value = "  -7.01%  "
$ws.Cells.Item(29, 5).Value = "  -6.85%  "
$ws.Cells.Item(30, 4).Value = "'8.88"
$ws.Cells.Item(30, 5).Value = "  -8.65%  "
$ws.Cells.Item(31, 4).Value = "'31.03"
$ws.Cells.Item(31, 5).Value = "  -7.98%  "
$ws.Cells.Item(32, 4).Value = "'3.92"
$ws.Cells.Item(32, 5).Value = "  -10.80%  "
$ws.Cells.Item(33, 5).Value = "  -7.83%  "
$ws.Cells.Item(34, 4).Value = "'11.23"
$ws.Cells.Item(34, 5).Value = "  -5.82%  "
$ws.Cells.Item(35, 4).Value = "'562.83"
$ws.Cells.Item(35, 5).Value = "  +12.26%  "
$ws.Cells.Item(36, 4).Value = "3.901.27"
$ws.Cells.Item(36, 5).Value = "  +2.62%  "
$ws.Cells.Item(37, 4).Value = "'0.107"
$ws.Cells.Item(37, 5).Value = "  -4.56%  "
$ws.Cells.Item(38, 4).Value = "'58.46"
$ws.Cells.Item(38, 5).Value = "  -6.32%  "
$ws.Cells.Item(39, 4).Value = "'1.00"
$ws.Cells.Item(39, 5).Value = "  -0.03%  "
$ws.Cells.Item(40, 4).Value = "'3.51"
$ws.Cells.Item(40, 5).Value = "  -5.38%  "
$ws.Cells.Item(41, 4).Value = "'3.54"
$ws.Cells.Item(41, 5).Value = "  +25.39%  "
$ws.Cells.Item(42, 5).Value = "  -11.28%  "
$ws.Cells.Item(43, 2).Value = "Kaspa"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43, 4).Value = "'0.129"
$ws.Cells.Item(43, 5).Value = "  -5.19%  "
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(44, 4).Value = "'2.71"
$ws.Cells.Item(44, 5).Value = "  -8.09%  "
$ws.Cells.Item(46, 4).Value = "'32.60"
$ws.Cells.Item(46, 5).Value = "  -6.24%  "
$ws.Cells.Item(47, 4).Value = "'0.0423"
$ws.Cells.Item(47, 5).Value = "  -7.92%  "
$ws.Cells.Item(48, 5).Value = "  -2.63%  "
$ws.Cells.Item(49, 5).Value = "  -8.32%  "
$ws.Cells.Item(50, 5).Value = "  -4.20%  "
$ws.Cells.Item(51, 4).Value = "'1.00"
$ws.Cells.Item(51, 5).Value = "  -0.15%  "
